# "Fruta / hortaliza, semanal"
#
# Inserts one new daily price record for Brócoli (Agrícola del Norte S.A.
# de Arica) as row 382, pushing the existing rows 382-405 down to 383-406.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 382 (shifts 382..405 -> 383..406).
$ws.Rows.Item(382).Insert()

# Populate the new row with the new record's data.
$ws.Cells.Item(382, 1).Value  = 1
$ws.Cells.Item(382, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(382, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(382, 4).Value  = 44746
$ws.Cells.Item(382, 5).Value  = 15
$ws.Cells.Item(382, 6).Value  = 100112023
$ws.Cells.Item(382, 7).Value  = "Brócoli"
$ws.Cells.Item(382, 8).Value  = "Sin especificar"
$ws.Cells.Item(382, 9).Value  = "Tercera"
$ws.Cells.Item(382, 10).Value = 1200
$ws.Cells.Item(382, 11).Value = 400
$ws.Cells.Item(382, 12).Value = 500
$ws.Cells.Item(382, 13).Value = 450
$ws.Cells.Item(382, 14).Value = "$/unidad"
$ws.Cells.Item(382, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(382, 16).Value = 450
$ws.Cells.Item(382, 17).Value = 1
$ws.Cells.Item(382, 18).Value = "Hortaliza"
